$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the data (row 2), pushing existing rows down
$ws.Range("A2:D5").EntireRow.Insert()

# The insert operation copies formatting from the row below; clear it so the
# new rows match the unstyled look of the rest of the data rows.
$ws.Range("A2:D5").ClearFormats()

$newRows = @(
    @("G1685AR.png", "1BEc5fw2mn1pWdTFtKCc1Ac5IPrCo3PFq"),
    @("G12301AR-.png", "1_UBXdXmyrPs0EwBnhV90mfhI3ZdSyn8H"),
    @("ISAALAMF16.png", "14BogV6rHDIcL1bJzvKgItflTszBh2Yjx"),
    @("EA5310MT.png", "168x1fBgYxstgD60EyOUekzYk21Cswij6")
)

$r = 2
foreach ($row in $newRows) {
    $name = $row[0]
    $id = $row[1]
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $id
    $ws.Cells.Item($r, 3).Value = "https://drive.google.com/file/d/$id/view?usp=drivesdk"
    $ws.Cells.Item($r, 4).Value = "https://drive.google.com/uc?export=view&id=$id"
    $r = $r + 1
}
